# Updated cryptos list with GitHub Actions - refresh prices / volume(1h)
# columns, plus swap the Aptos / Binance-PegBSC-USD rows (35 & 36).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.906.76"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "3.742.93"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'593.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "'166.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("D7").Value = "3.742.70"
$ws.Range("E7").Value = "  -2.11%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("E10").Value = "  -3.70%  "

$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").Value = "'0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.33%  "

$ws.Range("E13").Value = "  -4.77%  "

$ws.Range("D14").Value = "'36.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "4.372.63"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").Value = "3.752.84"
$ws.Range("E16").Value = "  -2.28%  "

$ws.Range("D17").Value = "67.900.47"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").Value = "'18.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "

$ws.Range("D19").Value = "'7.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.62%  "

$ws.Range("D21").Value = "'10.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("D22").Value = "'467.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  -5.17%  "

$ws.Range("D24").Value = "'83.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("E25").Value = "  -3.99%  "

$ws.Range("D26").Value = "'0.0000135"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.53%  "

$ws.Range("D27").Value = "'12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "'10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "3.891.89"
$ws.Range("E30").Value = "  -2.03%  "

$ws.Range("E31").Value = "  -4.89%  "

$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("E33").Value = "  -5.38%  "

$ws.Range("D34").Value = "'29.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'9.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  --%  "

$ws.Range("D37").Value = "3.696.82"
$ws.Range("E37").Value = "  -2.49%  "

$ws.Range("E38").Value = "  -4.50%  "

$ws.Range("D39").Value = "'3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.72%  "

$ws.Range("E40").Value = "  -1.17%  "

$ws.Range("D41").Value = "'0.993"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.30%  "

$ws.Range("D42").Value = "'5.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("D46").Value = "'8.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("D48").Value = "'45.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.74%  "

$ws.Range("D49").Value = "'393.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.34%  "

$ws.Range("D50").Value = "'143.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "'25.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
